# Fixed update to excel issue
#
# 1. Rename the "Requested quantity" headers on the existing sheets.
# 2. Add a new "PO Forecast" sheet (at the end) with the PO forecast data.

$wb = $excel.ActiveWorkbook

# --- 1. Rename existing headers -------------------------------------------
$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Capture the existing header / date formatting *before* adding the new
# sheet -- worksheet references read afresh after Worksheets.Add() can
# report stale/default formatting for other sheets.
$headerHAlign  = $wsWeekly.Range("B1").HorizontalAlignment
$headerVAlign  = $wsWeekly.Range("B1").VerticalAlignment
$dateNumberFmt = $wsWeekly.Range("A2").NumberFormat

# --- 2. Add the "PO Forecast" sheet ----------------------------------------
$wsForecast = $wb.Worksheets.Add()
$wsForecast.Name = "PO Forecast"

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$forecastData = @(
    @(45445.99999999999, 236, 165.4119923190973, 303.3076820452268),
    @(45508.99999999999, 153, 82.08832356171898, 224.2520079927008),
    @(45515.99999999999, 144, 77.00439212818584, 215.0087182320688),
    @(45578.99999999999, 61, -8.780719000329084, 134.2511777463187),
    @(45585.99999999999, 52, -18.17679923123305, 124.7265385945895),
    @(45599.99999999999, 33, -38.30827204519854, 102.7598390994168),
    @(45606.99999999999, 24, -43.21603189722975, 97.53984980463623),
    @(45613.99999999999, 15, -55.48734328240284, 84.54267306706453),
    @(45620.99999999999, 6, -65.0933717912435, 79.32069324275334),
    @(45627.99999999999, 0, -75.45007976725827, 65.15159799493622),
    @(45634.99999999999, 0, -77.08543811667437, 54.52832337476761),
    @(45641.99999999999, 0, -92.32890042724217, 46.58184636936019),
    @(45648.99999999999, 0, -103.3215383459327, 36.9387129846398),
    @(45655.99999999999, 0, -112.4078564027689, 30.43065616343274),
    @(45662.99999999999, 0, -120.7287656697358, 23.88639436142696),
    @(45669.99999999999, 0, -131.8369033066456, 13.15776412029337),
    @(45676.99999999999, 0, -137.7184540786214, 1.165188836497326)
)

$row = 2
foreach ($item in $forecastData) {
    $wsForecast.Cells.Item($row, 1).Value = $item[0]
    $wsForecast.Cells.Item($row, 2).Value = $item[1]
    $wsForecast.Cells.Item($row, 3).Value = $item[2]
    $wsForecast.Cells.Item($row, 4).Value = $item[3]
    $row++
}

# Match the header formatting used on the other sheets (bold, centered,
# thin border) and the date format used for the "ds" / date columns.
$wsForecast.Range("A1:D1").Font.Bold = $true
$wsForecast.Range("A1:D1").HorizontalAlignment = $headerHAlign
$wsForecast.Range("A1:D1").VerticalAlignment = $headerVAlign
$wsForecast.Range("A1:D1").Borders.LineStyle = 1

$wsForecast.Range("A2:A18").NumberFormat = $dateNumberFmt

# Move the new sheet to the end (after "Monthly Trend"), matching the
# workbook's sheet order in the target file.
$wsForecast.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

Write-Host "PO Forecast sheet added and headers renamed."
